# Update the COFINA Gabon representative's civility, name and title in the
# "ENTRE LES SOUSSIGNES" clause:
#   " est représentée par Monsieur El Hadji Mamadou FAYE, son Directeur Général, "
# becomes
#   " est représentée par Madame Jenny MVOU, son Directeur Général Adjointe, "

$d = $word.ActiveDocument

# 1) Drop the "Monsieur" civility from the (non-bold) lead-in text.
$d.Content.Find.Execute(" est représentée par Monsieur ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " est représentée par ", 2)

# 2) Insert the new civility "Madame " right before the bold name run so it
#    inherits the surrounding (non-bold) formatting.
$rng = $d.Content
$rng.Find.Execute("El Hadji Mamadou FAYE", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
  $rng.InsertBefore("Madame ")
}

# 3) Replace the old name with the new one, keeping the bold formatting of
#    the run that held the name.
$d.Content.Find.Execute("El Hadji Mamadou FAYE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jenny MVOU", 2)

# 4) Update the job title from "Directeur Général" to "Directeur Général Adjointe".
$d.Content.Find.Execute(", son Directeur Général, ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", son Directeur Général Adjointe, ", 2)
